$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "95.186.79"
Set-TextValue "D3" "3.590.52"
Set-TextValue "E3" "  +7.02%  "
Set-TextValue "E4" "  -0.11%  "
Set-TextValue "D5" "239.85"
Set-TextValue "E5" "  +3.50%  "
Set-TextValue "D6" "649.35"
Set-TextValue "E6" "  +5.28%  "
Set-TextValue "E7" "  +6.68%  "
Set-TextValue "D8" "0.407"
Set-TextValue "E8" "  +5.28%  "
Set-TextValue "D9" "0.999"
Set-TextValue "E9" "  -0.18%  "
Set-TextValue "D10" "0.998"
Set-TextValue "E10" "  +5.50%  "
Set-TextValue "D11" "3.586.23"
Set-TextValue "E11" "  +6.94%  "
Set-TextValue "D12" "42.94"
Set-TextValue "E12" "  +1.05%  "
Set-TextValue "D13" "0.199"
Set-TextValue "E13" "  +1.53%  "
Set-TextValue "D14" "6.31"
Set-TextValue "E14" "  +2.03%  "
Set-TextValue "D15" "4.268.20"
Set-TextValue "E15" "  +7.26%  "
Set-TextValue "D16" "94.923.80"
Set-TextValue "E16" "  +2.57%  "
Set-TextValue "E17" "  +4.88%  "
Set-TextValue "D18" "3.595.06"
Set-TextValue "E18" "  +7.02%  "
Set-TextValue "E19" "  -1.40%  "
Set-TextValue "D20" "12.48"
Set-TextValue "E20" "  +11.45%  "
Set-TextValue "D21" "18.02"
Set-TextValue "E21" "  +4.06%  "
Set-TextValue "D22" "3.47"
Set-TextValue "E22" "  +4.23%  "
Set-TextValue "E23" "  +13.56%  "
Set-TextValue "D24" "510.06"
Set-TextValue "E24" "  +3.24%  "
Set-TextValue "D25" "0.0000195"
Set-TextValue "E25" "  +7.40%  "
Set-TextValue "E26" "  +1.15%  "
Set-TextValue "D27" "92.77"
Set-TextValue "E27" "  +0.12%  "
Set-TextValue "D28" "12.80"
Set-TextValue "E28" "  +7.06%  "
Set-TextValue "D29" "3.10"
Set-TextValue "E29" "  +16.26%  "
Set-TextValue "D30" "11.28"
Set-TextValue "E30" "  +2.38%  "
Set-TextValue "D31" "0.999"
Set-TextValue "E31" "  -0.07%  "
Set-TextValue "E32" "  +2.32%  "
Set-TextValue "D33" "0.996"
Set-TextValue "E33" "  -1.01%  "
Set-TextValue "D34" "0.175"
Set-TextValue "E34" "  +2.04%  "
Set-TextValue "D35" "31.66"
Set-TextValue "E35" "  +11.15%  "
Set-TextValue "D36" "0.557"
Set-TextValue "E36" "  +6.52%  "
Set-TextValue "D37" "8.16"
Set-TextValue "E37" "  +9.94%  "
Set-TextValue "D38" "558.60"
Set-TextValue "E38" "  +1.13%  "
Set-TextValue "E39" "  +6.34%  "
Set-TextValue "B41" "ARBITRUM"
Set-TextValue "C41" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D41" "0.927"
Set-TextValue "E41" "  +5.66%  "
Set-TextValue "B42" "Kaspa"
Set-TextValue "C42" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D42" "0.150"
Set-TextValue "E42" "  +0.96%  "
Set-TextValue "E43" "  +2.45%  "
Set-TextValue "D44" "23.72"
Set-TextValue "E44" "  +0.28%  "
Set-TextValue "D45" "5.68"
Set-TextValue "E45" "  +5.90%  "
Set-TextValue "B46" "Stacks"
Set-TextValue "C46" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D46" "2.26"
Set-TextValue "E46" "  +8.27%  "
Set-TextValue "B47" "VeChain"
Set-TextValue "C47" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D47" "0.0417"
Set-TextValue "E47" "  +3.10%  "
Set-TextValue "D48" "54.34"
Set-TextValue "E48" "  +3.45%  "
Set-TextValue "D49" "32.77"
Set-TextValue "E49" "  +42.50%  "
Set-TextValue "D50" "3.44"
Set-TextValue "E50" "  -3.63%  "
Set-TextValue "D51" "8.07"
Set-TextValue "E51" "  +2.79%  "
